$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.653906
$ws.Range("H2").Value = 1.961718
$ws.Range("I2").Value = 0.00670030715761011
$ws.Range("J2").Value = 0.00670030715761011
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 46.33695966666667
$ws.Range("N2").Value = 139.010879
$ws.Range("O2").Value = 0.1993490803952133
$ws.Range("P2").Value = 0.1993490803952133
$ws.Range("Q2").Value = 30.30001594779133
$ws.Range("R2").Value = 272.700143530122
$ws.Range("S2").Value = 0.001335700070235041
$ws.Range("T2").Value = 0.001335700070235041

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.653906
$ws.Range("H3").Value = 1.961718
$ws.Range("I3").Value = 0.00670030715761011
$ws.Range("J3").Value = 0.00670030715761011
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 84.50960033333332
$ws.Range("N3").Value = 253.528801
$ws.Range("O3").Value = 0.3635739425333109
$ws.Range("P3").Value = 0.3635739425333109
$ws.Range("Q3").Value = 55.26133471556866
$ws.Range("R3").Value = 497.352012440118
$ws.Range("S3").Value = 0.00243605708947647
$ws.Range("T3").Value = 0.00243605708947647

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.653906
$ws.Range("H4").Value = 1.961718
$ws.Range("I4").Value = 0.00670030715761011
$ws.Range("J4").Value = 0.00670030715761011
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 72.52790466666666
$ws.Range("N4").Value = 217.583714
$ws.Range("O4").Value = 0.3120267536390091
$ws.Range("P4").Value = 0.3120267536390091
$ws.Range("Q4").Value = 47.42643202896133
$ws.Range("R4").Value = 426.8378882606519
$ws.Range("S4").Value = 0.002090675090773299
$ws.Range("T4").Value = 0.002090675090773299

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.653906
$ws.Range("H5").Value = 1.961718
$ws.Range("I5").Value = 0.00670030715761011
$ws.Range("J5").Value = 0.00670030715761011
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 29.06683666666666
$ws.Range("N5").Value = 87.20050999999998
$ws.Range("O5").Value = 0.1250502234324667
$ws.Range("P5").Value = 0.1250502234324667
$ws.Range("Q5").Value = 19.00697889735333
$ws.Range("R5").Value = 171.0628100761799
$ws.Range("S5").Value = 0.0008378749071253004
$ws.Range("T5").Value = 0.0008378749071253004

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 96.11977900000001
$ws.Range("H6").Value = 288.359337
$ws.Range("I6").Value = 0.984900036429704
$ws.Range("J6").Value = 0.984900036429704
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 46.33695966666667
$ws.Range("N6").Value = 139.010879
$ws.Range("O6").Value = 0.1993490803952133
$ws.Range("P6").Value = 0.1993490803952133
$ws.Range("Q6").Value = 4453.898322691914
$ws.Range("R6").Value = 40085.08490422722
$ws.Range("S6").Value = 0.1963389165434736
$ws.Range("T6").Value = 0.1963389165434736

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 96.11977900000001
$ws.Range("H7").Value = 288.359337
$ws.Range("I7").Value = 0.984900036429704
$ws.Range("J7").Value = 0.984900036429704
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 84.50960033333332
$ws.Range("N7").Value = 253.528801
$ws.Range("O7").Value = 0.3635739425333109
$ws.Range("P7").Value = 0.3635739425333109
$ws.Range("Q7").Value = 8123.044107418326
$ws.Range("R7").Value = 73107.39696676494
$ws.Range("S7").Value = 0.358083989245949
$ws.Range("T7").Value = 0.358083989245949

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 96.11977900000001
$ws.Range("H8").Value = 288.359337
$ws.Range("I8").Value = 0.984900036429704
$ws.Range("J8").Value = 0.984900036429704
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 72.52790466666666
$ws.Range("N8").Value = 217.583714
$ws.Range("O8").Value = 0.3120267536390091
$ws.Range("P8").Value = 0.3120267536390091
$ws.Range("Q8").Value = 6971.366167893068
$ws.Range("R8").Value = 62742.29551103762
$ws.Range("S8").Value = 0.3073151610261023
$ws.Range("T8").Value = 0.3073151610261023

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 96.11977900000001
$ws.Range("H9").Value = 288.359337
$ws.Range("I9").Value = 0.984900036429704
$ws.Range("J9").Value = 0.984900036429704
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 29.06683666666666
$ws.Range("N9").Value = 87.20050999999998
$ws.Range("O9").Value = 0.1250502234324667
$ws.Range("P9").Value = 0.1250502234324667
$ws.Range("Q9").Value = 2793.897916629096
$ws.Range("R9").Value = 25145.08124966187
$ws.Range("S9").Value = 0.1231619696141791
$ws.Range("T9").Value = 0.1231619696141791

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.2341223333333333
$ws.Range("H10").Value = 0.702367
$ws.Range("I10").Value = 0.002398955730318598
$ws.Range("J10").Value = 0.002398955730318598
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 46.33695966666667
$ws.Range("N10").Value = 139.010879
$ws.Range("O10").Value = 0.1993490803952133
$ws.Range("P10").Value = 0.1993490803952133
$ws.Range("Q10").Value = 10.84851711673256
$ws.Range("R10").Value = 97.63665405059298
$ws.Range("S10").Value = 0.00047822961874784
$ws.Range("T10").Value = 0.0004782296187478398

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.2341223333333333
$ws.Range("H11").Value = 0.702367
$ws.Range("I11").Value = 0.002398955730318598
$ws.Range("J11").Value = 0.002398955730318598
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 84.50960033333332
$ws.Range("N11").Value = 253.528801
$ws.Range("O11").Value = 0.3635739425333109
$ws.Range("P11").Value = 0.3635739425333109
$ws.Range("Q11").Value = 19.78558481910744
$ws.Range("R11").Value = 178.070263371967
$ws.Range("S11").Value = 0.0008721977928348109
$ws.Range("T11").Value = 0.0008721977928348108

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.2341223333333333
$ws.Range("H12").Value = 0.702367
$ws.Range("I12").Value = 0.002398955730318598
$ws.Range("J12").Value = 0.002398955730318598
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 72.52790466666666
$ws.Range("N12").Value = 217.583714
$ws.Range("O12").Value = 0.3120267536390091
$ws.Range("P12").Value = 0.3120267536390091
$ws.Range("Q12").Value = 16.98040227233755
$ws.Range("R12").Value = 152.823620451038
$ws.Range("S12").Value = 0.0007485383686550104
$ws.Range("T12").Value = 0.0007485383686550103

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.2341223333333333
$ws.Range("H13").Value = 0.702367
$ws.Range("I13").Value = 0.002398955730318598
$ws.Range("J13").Value = 0.002398955730318598
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 29.06683666666666
$ws.Range("N13").Value = 87.20050999999998
$ws.Range("O13").Value = 0.1250502234324667
$ws.Range("P13").Value = 0.1250502234324667
$ws.Range("Q13").Value = 6.805195623018887
$ws.Range("R13").Value = 61.24676060716998
$ws.Range("S13").Value = 0.0002999899500809372
$ws.Range("T13").Value = 0.0002999899500809371

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 0.5856290000000001
$ws.Range("H14").Value = 1.756887
$ws.Range("I14").Value = 0.00600070068236727
$ws.Range("J14").Value = 0.006000700682367269
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 46.33695966666667
$ws.Range("N14").Value = 139.010879
$ws.Range("O14").Value = 0.1993490803952133
$ws.Range("P14").Value = 0.1993490803952133
$ws.Range("Q14").Value = 27.13626735263033
$ws.Range("R14").Value = 244.226406173673
$ws.Range("S14").Value = 0.001196234162756844
$ws.Range("T14").Value = 0.001196234162756844

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 0.5856290000000001
$ws.Range("H15").Value = 1.756887
$ws.Range("I15").Value = 0.00600070068236727
$ws.Range("J15").Value = 0.006000700682367269
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 84.50960033333332
$ws.Range("N15").Value = 253.528801
$ws.Range("O15").Value = 0.3635739425333109
$ws.Range("P15").Value = 0.3635739425333109
$ws.Range("Q15").Value = 49.49127273360967
$ws.Range("R15").Value = 445.421454602487
$ws.Range("S15").Value = 0.002181698405050597
$ws.Range("T15").Value = 0.002181698405050597

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 0.5856290000000001
$ws.Range("H16").Value = 1.756887
$ws.Range("I16").Value = 0.00600070068236727
$ws.Range("J16").Value = 0.006000700682367269
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 72.52790466666666
$ws.Range("N16").Value = 217.583714
$ws.Range("O16").Value = 0.3120267536390091
$ws.Range("P16").Value = 0.3120267536390091
$ws.Range("Q16").Value = 42.47444428203534
$ws.Range("R16").Value = 382.269998538318
$ws.Range("S16").Value = 0.001872379153478446
$ws.Range("T16").Value = 0.001872379153478445

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 0.5856290000000001
$ws.Range("H17").Value = 1.756887
$ws.Range("I17").Value = 0.00600070068236727
$ws.Range("J17").Value = 0.006000700682367269
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 29.06683666666666
$ws.Range("N17").Value = 87.20050999999998
$ws.Range("O17").Value = 0.1250502234324667
$ws.Range("P17").Value = 0.1250502234324667
$ws.Range("Q17").Value = 17.02238249026333
$ws.Range("R17").Value = 153.20144241237
$ws.Range("S17").Value = 0.0007503889610813827
$ws.Range("T17").Value = 0.0007503889610813826
